$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2673.9
$ws.Range("G2").Value = 300

$ws.Range("F3").Value = 5289
$ws.Range("G3").Value = 474

$ws.Range("F4").Value = 2067
$ws.Range("G4").Value = 168

$ws.Range("F5").Value = 5096.4
$ws.Range("G5").Value = 420

$ws.Range("F6").Value = 4837.92
$ws.Range("G6").Value = 462

$ws.Range("F7").Value = 13386.18
$ws.Range("G7").Value = 924

$ws.Range("F8").Value = 4216.56
$ws.Range("G8").Value = 474

$ws.Range("F9").Value = 218.4
$ws.Range("G9").Value = 12

$ws.Range("F10").Value = 28.2
$ws.Range("G10").Value = 6

$ws.Range("F11").Value = 10858.8
$ws.Range("G11").Value = 924
